# Update Datasheet / Supplier links for the oscillator rows (X1, X2)
# from Abracon parts to the new Kyocera AVX parts, per upstream BoM change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# Row 26 -> "18" / ASCO-Oscillator / X1 / 22.5792M
$ws.Range("I26").Value = "https://media.digikey.com/pdf/Data%20Sheets/Kyocera%20International/Z_Series_X_Type.pdf"
$ws.Range("J26").Value = "https://www.digikey.ch/en/products/detail/kyocera-avx/KC2016Z22-5792C1KX00/11610237"

# Row 27 -> "19" / X2 / 24.576M
$ws.Range("I27").Value = "https://media.digikey.com/pdf/Data%20Sheets/Kyocera%20International/Z_Series_X_Type.pdf"
$ws.Range("J27").Value = "https://www.digikey.ch/en/products/detail/kyocera-avx/KC2016Z24-5760C1KX00/11610181"
